$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:B3").Copy()
$ws.Range("Z1:Z1").PasteSpecial(-4122)

$ws.Rows("3:7").Delete()

$ws.Range("C2:I2").Copy()
$ws.Range("C3:I3").PasteSpecial(-4122)

$ws.Range("Z1:Z1").Copy()
$ws.Range("B3:B3").PasteSpecial(-4122)
$ws.Range("B3:B3").ClearContents()

$ws.Range("Z1:Z1").Clear()

$ws.Range("C4:I4").Value = 0

$null = $ws.Range("I4").Select()
